$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.467.17"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "'3.417.74"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'563.64"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").Value = "'175.54"
$ws.Range("E6").Value = "  +3.54%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +3.62%  "
$ws.Range("D8").Value = "'3.411.59"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +13.01%  "
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("D12").Value = "'54.98"
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").Value = "'0.0000280"
$ws.Range("E13").Value = "  +6.68%  "
$ws.Range("D14").Value = "'9.18"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("D15").Value = "'3.948.35"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "'18.43"
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("D17").Value = "'3.401.25"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("D19").Value = "'11.97"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").Value = "'65.403.37"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("D22").Value = "'470.42"
$ws.Range("E22").Value = "  +16.77%  "
$ws.Range("D23").Value = "'4.99"
$ws.Range("E23").Value = "  +17.73%  "
$ws.Range("D24").Value = "'4.16"
$ws.Range("E24").Value = "  +3.51%  "
$ws.Range("D25").Value = "'86.69"
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("D26").Value = "'13.81"
$ws.Range("E26").Value = "  +5.49%  "
$ws.Range("E27").Value = "  +4.09%  "
$ws.Range("D28").Value = "'2.91"
$ws.Range("E28").Value = "  +7.86%  "
$ws.Range("D29").Value = "'8.91"
$ws.Range("E29").Value = "  +4.67%  "
$ws.Range("D30").Value = "'30.81"
$ws.Range("E30").Value = "  +6.67%  "
$ws.Range("D31").Value = "'6.74"
$ws.Range("E31").Value = "  +5.34%  "
$ws.Range("D32").Value = "'11.59"
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("D33").Value = "'586.63"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D35").Value = "'60.40"
$ws.Range("E35").Value = "  +5.54%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'36.18"
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.49"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").Value = "'0.0₃0756"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").Value = "'0.377"
$ws.Range("E41").Value = "  +3.95%  "
$ws.Range("D42").Value = "'3.117.74"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("E45").Value = "  +4.82%  "
$ws.Range("D46").Value = "'0.0417"
$ws.Range("E46").Value = "  +4.75%  "
$ws.Range("D47").Value = "'3.23"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("E48").Value = "  +6.03%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'137.02"
$ws.Range("E50").Value = "  +3.48%  "
$ws.Range("E51").Value = "  +5.47%  "
